$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 1.27
$ws.Range("H4").Value = 4.8
$ws.Range("I4").Value = 12
$ws.Range("J4").Value = 1.75
$ws.Range("K4").Value = 2.27
$ws.Range("L4").Value = 10.25
$ws.Range("N4").Value = 7.2
$ws.Range("O4").Value = 1.31
$ws.Range("P4").Value = 3.15
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.82
$ws.Range("S4").Value = 1.42
$ws.Range("T4").Value = 2.65
$ws.Range("U4").Value = 2.65
$ws.Range("V4").Value = 1.42
$ws.Range("X4").Value = 5
$ws.Range("Y4").Value = 9.75
$ws.Range("Z4").Value = 6.8
$ws.Range("AA4").Value = 13
$ws.Range("AC4").Value = 7.2
$ws.Range("AD4").Value = 10.25
$ws.Range("AE4").Value = 37
$ws.Range("AF4").Value = 250
$ws.Range("AH4").Value = 22
$ws.Range("AI4").Value = 90
$ws.Range("AJ4").Value = 40
$ws.Range("AK4").Value = 500
$ws.Range("AM4").Value = 200
$ws.Range("AN4").Value = 2.82
$ws.Range("AO4").Value = 5.5
$ws.Range("AQ4").Value = 15.5
$ws.Range("AT4").Value = 2.65
$ws.Range("AU4").Value = 11.25
$ws.Range("AW4").Value = 11.25
$ws.Range("AX4").Value = 90
$ws.Range("AY4").Value = 90

# Row 6 updates
$ws.Range("N6").Value = 13
$ws.Range("O6").Value = 1.18
$ws.Range("P6").Value = 4.5
$ws.Range("S6").Value = 1.3
$ws.Range("T6").Value = 3.4
$ws.Range("U6").Value = 1.93
$ws.Range("V6").Value = 1.82

# Row 12 updates
$ws.Range("Q12").Value = 1.9
$ws.Range("R12").Value = 1.9
